$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, pushing existing rows 19-26 down to 20-27.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new data record.
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 45258
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100103
$ws.Cells.Item(19, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(19, 9).Value = 100103001
$ws.Cells.Item(19, 10).Value = "Cereza"
$ws.Cells.Item(19, 11).Value = "Lapins"
$ws.Cells.Item(19, 12).Value = "Segunda"
$ws.Cells.Item(19, 13).Value = 330
$ws.Cells.Item(19, 14).Value = 14000
$ws.Cells.Item(19, 15).Value = 15000
$ws.Cells.Item(19, 16).Value = 14455
$ws.Cells.Item(19, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 1446
$ws.Cells.Item(19, 20).Value = 10
